# The presentation ships two DrawingML themes:
#   ppt/theme/theme1.xml  -> "Office Theme" palette (used by the Notes Master)
#   ppt/theme/theme2.xml  -> "Integral"     palette (used by the Slide Master
#                             / the deck's design, i.e. what Slide.ThemeColorScheme
#                             and friends operate on)
#
# The authored edit swaps the two themes' contents, so the deck's design
# (theme2.xml) ends up with the colors that used to live in theme1.xml
# ("Office Theme" palette), and vice-versa. Apply that swap to the design
# theme through the PowerPoint color-scheme object model: each of the 12
# theme color slots is addressed in the fixed order
#   1 dk1, 2 lt1, 3 dk2, 4 lt2,
#   5 accent1, 6 accent2, 7 accent3, 8 accent4, 9 accent5, 10 accent6,
#   11 hlink, 12 folHlink
# and re-pointed to the RGB value the "Office Theme" palette used.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

$tcs.Item(1).RGB = 0
$tcs.Item(2).RGB = 16777215
$tcs.Item(3).RGB = 6968388
$tcs.Item(4).RGB = 15132391
$tcs.Item(5).RGB = 13998939
$tcs.Item(6).RGB = 3243501
$tcs.Item(7).RGB = 10855845
$tcs.Item(8).RGB = 49407
$tcs.Item(9).RGB = 12874308
$tcs.Item(10).RGB = 4697456
$tcs.Item(11).RGB = 12673797
$tcs.Item(12).RGB = 7491477
